$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row additions (row 1), matching style of existing header cells (e.g. H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header style from an existing header cell (H1) to the new header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Data row additions (row 2)
$ws.Range("I2").Value = 9
$ws.Range("J2").Value = 9
